$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.183.36"
$ws.Range("E2").Value = "  -0.34%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.642.48"
$ws.Range("E3").Value = "  +0.07%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.52"
$ws.Range("E5").Value = "  -0.65%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.02"
$ws.Range("E6").Value = "  +3.05%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("E8").Value = "  -0.77%  "

$ws.Range("E9").Value = "  -1.46%  "

$ws.Range("E10").Value = "  -1.49%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.27"
$ws.Range("E11").Value = "  +0.19%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.351"
$ws.Range("E12").Value = "  -1.07%  "

$ws.Range("E13").Value = "  -1.05%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.122.93"
$ws.Range("E14").Value = "  +0.06%  "

$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000188"
$ws.Range("E15").Value = "  -2.33%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "68.045.15"
$ws.Range("E16").Value = "  -0.35%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.660.09"
$ws.Range("E17").Value = "  +0.61%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.37"
$ws.Range("E18").Value = "  -1.36%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "360.73"
$ws.Range("E19").Value = "  -1.63%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.34"
$ws.Range("E20").Value = "  -2.49%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.38"
$ws.Range("E21").Value = "  -0.91%  "

$ws.Range("E22").Value = "  -2.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.07"
$ws.Range("E23").Value = "  -0.69%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.78"
$ws.Range("E24").Value = "  +1.10%  "

$ws.Range("E25").Value = "  +0.15%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.78"
$ws.Range("E26").Value = "  -1.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.777.00"
$ws.Range("E27").Value = "  +0.17%  "

$ws.Range("E28").Value = "  -2.37%  "

$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.993"
$ws.Range("E29").Value = "  -0.88%  "

$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "564.70"
$ws.Range("E30").Value = "  -2.22%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.97"
$ws.Range("E31").Value = "  -1.68%  "

$ws.Range("E32").Value = "  -3.23%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.85"
$ws.Range("E33").Value = "  -1.55%  "

$ws.Range("E35").Value = "  -3.53%  "

$ws.Range("E36").Value = "  -3.82%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "160.47"
$ws.Range("E37").Value = "  -0.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.75"
$ws.Range("E38").Value = "  +1.36%  "

$ws.Range("E39").Value = "  -0.83%  "

$ws.Range("E40").Value = "  -1.52%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.30"
$ws.Range("E41").Value = "  -2.50%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.80"
$ws.Range("E42").Value = "  +0.41%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.61"
$ws.Range("E43").Value = "  -3.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₆0326"
$ws.Range("E44").Value = "  +0.14%  "

$ws.Range("E45").Value = "  +0.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "157.13"
$ws.Range("E46").Value = "  -0.44%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.74"
$ws.Range("E47").Value = "  -1.41%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.01"
$ws.Range("E48").Value = "  +0.05%  "

$ws.Range("E49").Value = "  -2.41%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0774"
$ws.Range("E50").Value = "  -1.68%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.614"
$ws.Range("E51").Value = "  -0.33%  "
